$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (row 1): append " (accuracy)" to each dataset/split header ---
$ws.Range("B1").Value = "Chinatown Train (accuracy)"
$ws.Range("C1").Value = "Chinatown Val (accuracy)"
$ws.Range("D1").Value = "Chinatown Test (accuracy)"
$ws.Range("E1").Value = "ECG200 Train (accuracy)"
$ws.Range("F1").Value = "ECG200 Val (accuracy)"
$ws.Range("G1").Value = "ECG200 Test (accuracy)"
$ws.Range("H1").Value = "ItalyPowerDemand Train (accuracy)"
$ws.Range("I1").Value = "ItalyPowerDemand Val (accuracy)"
$ws.Range("J1").Value = "ItalyPowerDemand Test (accuracy)"

# --- Rename model labels (use hyphen instead of underscore) ---
$ws.Range("A3").Value = "decision-tree"

# --- Add new row 4 for the linear-regression model results ---
$ws.Range("A4").Value = "linear-regression"
$ws.Range("B4").Value = 1
$ws.Range("C4").Value = 0.86274509803921495
$ws.Range("D4").Value = 0.89211618257261405
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.53333333333333299
$ws.Range("G4").Value = 0.28571428571428498
$ws.Range("H4").Value = 1
$ws.Range("I4").Value = 0.50324675324675305
$ws.Range("J4").Value = 0.50069348127600499

# --- Update the selected/active cell to match the recorded view state ---
$ws.Range("G10").Select()

$wb.Save()
